$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.00018693
$ws.Range("F2").Value = 0.01265751
$ws.Range("G2").Value = 0.00034098274285714286

$ws.Range("E3").Value = 0.00292806
$ws.Range("F3").Value = 0.00702603
$ws.Range("G3").Value = 0.0032815819513797636
